$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / rId1) - first worksheet in the workbook
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 139
$ws1.Range("F3").Value = 1324
$ws1.Range("F4").Value = 1123
$ws1.Range("F5").Value = 1014
$ws1.Range("F6").Value = 1790
$ws1.Range("F7").Value = 556
$ws1.Range("F8").Value = 1193
$ws1.Range("F12").Value = 292
$ws1.Range("F13").Value = 64
$ws1.Range("F15").Value = 685
$ws1.Range("F16").Value = 165
$ws1.Range("F21").Value = 143
$ws1.Range("F22").Value = 668
$ws1.Range("F23").Value = 34
$ws1.Range("F24").Value = 643
$ws1.Range("F27").Value = 872
$ws1.Range("F29").Value = 159

# Sheet "全部类型" (sheet4 / rId4) - contains the same events shifted by one row
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value = 139
$ws4.Range("F4").Value = 1324
$ws4.Range("F5").Value = 1123
$ws4.Range("F6").Value = 1014
$ws4.Range("F7").Value = 1790
$ws4.Range("F8").Value = 556
$ws4.Range("F9").Value = 1193
$ws4.Range("F14").Value = 292
$ws4.Range("F15").Value = 64
$ws4.Range("F17").Value = 685
$ws4.Range("F18").Value = 165
$ws4.Range("F29").Value = 143
$ws4.Range("F30").Value = 668
$ws4.Range("F31").Value = 34
$ws4.Range("F32").Value = 643
$ws4.Range("F35").Value = 872
$ws4.Range("F39").Value = 159
